$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(27, 'Metal Mario - Tennis', 'Metal Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(28, 'Super Mario Cereal', 'Mario Cereal', 'Others', 'Kellogs', 'Card'),
    @(29, 'Baby Mario - Soccer', 'Baby Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(30, 'Metal Mario - Soccer', 'Metal Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(31, 'Mario - Soccer', 'Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(32, 'Mario', 'Mario', 'Super Mario Bros.', 'Super Mario', 'Figure'),
    @(33, '8-Bit Mario Modern Color', 'Mario', '8-bit Mario', 'Super Mario', 'Figure'),
    @(34, 'Dr. Mario', 'Mario', 'Super Smash Bros.', 'Super Mario', 'Figure'),
    @(35, 'Baby Mario - Horse Racing', 'Baby Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(36, 'Mario - Cat', 'Mario', 'Super Mario Bros.', 'Super Mario', 'Figure'),
    @(37, 'Baby Mario - Golf', 'Baby Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(38, 'Mario - Wedding', 'Mario', 'Super Mario Bros.', 'Super Mario', 'Figure'),
    @(39, 'Metal Mario - Golf', 'Metal Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(40, 'Mario', 'Mario', 'Super Smash Bros.', 'Super Mario', 'Figure'),
    @(41, 'Mario - Tennis', 'Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(42, 'Baby Mario - Tennis', 'Baby Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(43, 'Mario - Gold Edition', 'Mario', 'Super Mario Bros.', 'Super Mario', 'Figure'),
    @(44, 'Mario - Power Up Band', 'Mario', 'Super Nintendo World', 'Super Mario', 'Band'),
    @(45, 'Metal Mario - Horse Racing', 'Metal Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(46, 'Baby Mario - Baseball', 'Baby Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(47, 'Mario - Golf', 'Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(48, 'Mario - Silver Edition', 'Mario', 'Super Mario Bros.', 'Super Mario', 'Figure'),
    @(49, 'Mario - Horse Racing', 'Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(50, '8-Bit Mario Classic Color', 'Mario', '8-bit Mario', 'Super Mario', 'Figure'),
    @(51, 'Metal Mario - Baseball', 'Metal Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
    @(52, 'Mario - Baseball', 'Mario', 'Mario Sports Superstars', 'Mario Sports Superstars', 'Card'),
)

$startRow = 28
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
